# "draft of w5 slides": fill in links for the w4 part-2 and w5 slide decks
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value  = "w4p2"
$ws.Range("D10").Value = "w5p1"
$ws.Range("D11").Value = "w5p2"

# Reflect where the author was looking/working when they saved: scrolled
# down a bit and with D12 selected in the frozen (right) pane.
$excel.Goto($ws.Range("D12"), $false)
